# Updated cryptos list on Fri Nov 29 10:48:23 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# cryptos table, and swaps the Cosmos / MantraDAO rows (50/51) to reflect
# their new ranking order.
#
# Price cells are stored as literal text in the source data (they use a
# "."-as-thousands-separator style, e.g. "96.790.33", and keep trailing
# zeros, e.g. "13.10"), so every D-column write below is given a leading
# apostrophe. That is the standard Excel "force text" entry convention:
# Excel strips the apostrophe and stores the remainder verbatim as a text
# value (quote-prefixed) instead of re-parsing it as a number - which
# preserves values like "13.10" / "0.999" / "96.790.33" exactly as typed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''96.790.33'
$ws.Range("E2").Value = '  +1.88%  '
$ws.Range("D3").Value = '''3.580.72'
$ws.Range("E3").Value = '  -0.83%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''241.54'
$ws.Range("E5").Value = '  +2.28%  '
$ws.Range("D6").Value = '''653.70'
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("E7").Value = '  +14.89%  '
$ws.Range("E8").Value = '  +2.43%  '
$ws.Range("E9").Value = '  +9.08%  '
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("D11").Value = '''3.578.26'
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("D12").Value = '''43.58'
$ws.Range("E12").Value = '  +2.90%  '
$ws.Range("E13").Value = '  +1.49%  '
$ws.Range("E14").Value = '  +1.56%  '
$ws.Range("D15").Value = '''4.243.15'
$ws.Range("E15").Value = '  -1.39%  '
$ws.Range("D16").Value = '''96.567.23'
$ws.Range("E16").Value = '  +1.70%  '
$ws.Range("E17").Value = '  +2.77%  '
$ws.Range("D18").Value = '''3.571.98'
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").Value = '''7.78'
$ws.Range("E19").Value = '  -1.58%  '
$ws.Range("D20").Value = '''12.64'
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("D21").Value = '''18.01'
$ws.Range("E21").Value = '  +0.57%  '
$ws.Range("D22").Value = '''0.545'
$ws.Range("E22").Value = '  +13.98%  '
$ws.Range("D23").Value = '''508.68'
$ws.Range("E23").Value = '  +0.88%  '
$ws.Range("D24").Value = '''3.40'
$ws.Range("E24").Value = '  -4.08%  '
$ws.Range("E25").Value = '  +5.80%  '
$ws.Range("D26").Value = '''0.0000201'
$ws.Range("E26").Value = '  +2.47%  '
$ws.Range("D27").Value = '''96.83'
$ws.Range("E27").Value = '  +1.64%  '
$ws.Range("D28").Value = '''13.10'
$ws.Range("E28").Value = '  +4.84%  '
$ws.Range("D29").Value = '''3.769.82'
$ws.Range("E29").Value = '  -0.99%  '
$ws.Range("D30").Value = '''0.154'
$ws.Range("E30").Value = '  +12.14%  '
$ws.Range("D31").Value = '''3.04'
$ws.Range("E31").Value = '  -2.76%  '
$ws.Range("D32").Value = '''11.51'
$ws.Range("E32").Value = '  +2.53%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  +4.21%  '
$ws.Range("D35").Value = '''0.999'
$ws.Range("D36").Value = '''31.42'
$ws.Range("E36").Value = '  -2.01%  '
$ws.Range("D37").Value = '''624.69'
$ws.Range("E37").Value = '  +10.08%  '
$ws.Range("D38").Value = '''8.87'
$ws.Range("E38").Value = '  +9.78%  '
$ws.Range("D39").Value = '''0.570'
$ws.Range("E39").Value = '  +2.59%  '
$ws.Range("E40").Value = '  +11.59%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("E42").Value = '  +1.47%  '
$ws.Range("D43").Value = '''0.908'
$ws.Range("E43").Value = '  -0.43%  '
$ws.Range("D44").Value = '''1.83'
$ws.Range("E44").Value = '  +6.12%  '
$ws.Range("D45").Value = '''5.79'
$ws.Range("E45").Value = '  +3.62%  '
$ws.Range("D46").Value = '''0.0431'
$ws.Range("E46").Value = '  +4.54%  '
$ws.Range("D47").Value = '''2.31'
$ws.Range("E47").Value = '  +4.16%  '
$ws.Range("D48").Value = '''23.55'
$ws.Range("E48").Value = '  -0.56%  '
$ws.Range("D49").Value = '''32.86'
$ws.Range("E49").Value = '  -7.15%  '

# Rows 50/51 swap places: MantraDAO <-> Cosmos, with refreshed price/volume.
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '''8.34'
$ws.Range("E50").Value = '  +4.59%  '
$ws.Range("B51").Value = 'MantraDAO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D51").Value = '''3.51'
$ws.Range("E51").Value = '  -0.88%  '
